$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2661
$ws.Range("I38").Value = 1041.5
$ws.Range("J38").Value = 5900
$ws.Range("K38").Value = 3124.5
$ws.Range("L38").Value = 17700
$ws.Range("M38").Value = -2752.5
$ws.Range("N38").Value = -18444
$ws.Range("H40").Value = 2800.4707
$ws.Range("J40").Value = 2948.6365
$ws.Range("L40").Value = 2948.6365
$ws.Range("N40").Value = -3298.6365
$ws.Range("H41").Value = 212.55556
$ws.Range("I41").Value = 135
$ws.Range("J41").Value = 309.5
$ws.Range("K41").Value = 135
$ws.Range("L41").Value = 309.5
$ws.Range("M41").Value = 305
$ws.Range("N41").Value = -1189.5
$ws.Range("H74").Value = 14855.875
$ws.Range("J74").Value = 9499
$ws.Range("L74").Value = 9499
$ws.Range("N74").Value = -11371
$ws.Range("H76").Value = 21502
$ws.Range("I76").Value = 20000
$ws.Range("J76").Value = 23004
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 23004
$ws.Range("M76").Value = -19685
$ws.Range("N76").Value = -23634
$ws.Range("H77").Value = 14855.875
$ws.Range("J77").Value = 9499
$ws.Range("L77").Value = 47495
$ws.Range("N77").Value = -56855
$ws.Range("H79").Value = 21502
$ws.Range("I79").Value = 20000
$ws.Range("J79").Value = 23004
$ws.Range("K79").Value = 20000
$ws.Range("L79").Value = 23004
$ws.Range("M79").Value = -18908
$ws.Range("N79").Value = -25188
$ws.Range("H94").Value = 621.9091
$ws.Range("I94").Value = 621.9091
$ws.Range("K94").Value = 621.9091
$ws.Range("M94").Value = -170.9091
$ws.Range("H137").Value = 9891.091
$ws.Range("I137").Value = 1879.7
$ws.Range("K137").Value = 5639.1
$ws.Range("M137").Value = -3089.1
$ws.Range("H138").Value = 5547.0586
$ws.Range("I138").Value = 3555.1428
$ws.Range("J138").Value = 6941.4
$ws.Range("K138").Value = 10665.4284
$ws.Range("L138").Value = 20824.2
$ws.Range("M138").Value = -5525.428400000001
$ws.Range("N138").Value = -31104.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1193823
$ws.Range("I32").Value = 589446.0600000001
$ws.Range("K32").Value = 589446.0600000001
$ws.Range("M32").Value = -589159.0600000001
$ws.Range("H61").Value = 3703.2273
$ws.Range("I61").Value = 2974.7273
$ws.Range("K61").Value = 2974.7273
$ws.Range("M61").Value = -2762.7273
$ws.Range("H63").Value = 2270.2856
$ws.Range("I63").Value = 2119.2
$ws.Range("J63").Value = 2648
$ws.Range("K63").Value = 2119.2
$ws.Range("L63").Value = 2648
$ws.Range("M63").Value = -1433.2
$ws.Range("N63").Value = -4020
$ws.Range("H66").Value = 2270.2856
$ws.Range("I66").Value = 2119.2
$ws.Range("J66").Value = 2648
$ws.Range("K66").Value = 10596
$ws.Range("L66").Value = 13240
$ws.Range("M66").Value = -7164
$ws.Range("N66").Value = -20104
$ws.Range("H136").Value = 3703.2273
$ws.Range("I136").Value = 2974.7273
$ws.Range("K136").Value = 8924.1819
$ws.Range("M136").Value = -6374.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2934.4075
$ws.Range("I86").Value = 2410.5789
$ws.Range("K86").Value = 2410.5789
$ws.Range("M86").Value = -1287.5789
$ws.Range("H89").Value = 2934.4075
$ws.Range("I89").Value = 2410.5789
$ws.Range("K89").Value = 12052.8945
$ws.Range("M89").Value = -6436.8945
$ws.Range("H94").Value = 33334374
$ws.Range("J94").Value = 998.5
$ws.Range("L94").Value = 998.5
$ws.Range("N94").Value = -1900.5
$ws.Range("H99").Value = 205180.8
$ws.Range("I99").Value = 254601
$ws.Range("K99").Value = 254601
$ws.Range("M99").Value = -253103
$ws.Range("H134").Value = 2254.1904
$ws.Range("I134").Value = 1641.6923
$ws.Range("K134").Value = 4925.0769
$ws.Range("M134").Value = -2390.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 547.1429000000001
$ws.Range("I7").Value = 670.5625
$ws.Range("K7").Value = 670.5625
$ws.Range("M7").Value = -557.5625
$ws.Range("H31").Value = 4170463.2
$ws.Range("I31").Value = 2649.5
$ws.Range("J31").Value = 4811665
$ws.Range("K31").Value = 2649.5
$ws.Range("L31").Value = 4811665
$ws.Range("M31").Value = -2354.5
$ws.Range("N31").Value = -4812255
$ws.Range("H34").Value = 4170463.2
$ws.Range("I34").Value = 2649.5
$ws.Range("J34").Value = 4811665
$ws.Range("K34").Value = 2649.5
$ws.Range("L34").Value = 4811665
$ws.Range("M34").Value = -2447.5
$ws.Range("N34").Value = -4812069
$ws.Range("H132").Value = 2980.4614
$ws.Range("I132").Value = 2291.111
$ws.Range("K132").Value = 6873.333
$ws.Range("M132").Value = -4343.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 236.25
$ws.Range("J50").Value = 236.25
$ws.Range("L50").Value = 708.75
$ws.Range("N50").Value = -1670.75
$ws.Range("H53").Value = 236.25
$ws.Range("J53").Value = 236.25
$ws.Range("L53").Value = 708.75
$ws.Range("N53").Value = -1670.75
$ws.Range("H68").Value = 1726668.4
$ws.Range("J68").Value = 1963430.8
$ws.Range("L68").Value = 5890292.4
$ws.Range("N68").Value = -5891914.4
$ws.Range("H71").Value = 1726668.4
$ws.Range("J71").Value = 1963430.8
$ws.Range("L71").Value = 17670877.2
$ws.Range("N71").Value = -17678989.2
$ws.Range("H107").Value = 2095.762
$ws.Range("J107").Value = 3149.1667
$ws.Range("L107").Value = 9447.500100000001
$ws.Range("N107").Value = -13287.5001
$ws.Range("H131").Value = 1554923.8
$ws.Range("I131").Value = 21144.7
$ws.Range("J131").Value = 2102702
$ws.Range("K131").Value = 63434.10000000001
$ws.Range("L131").Value = 6308106
$ws.Range("M131").Value = -58394.10000000001
$ws.Range("N131").Value = -6318186
$ws.Range("H132").Value = 5542.1465
$ws.Range("I132").Value = 5867.625
$ws.Range("J132").Value = 5463.242
$ws.Range("K132").Value = 52808.625
$ws.Range("L132").Value = 49169.178
$ws.Range("M132").Value = -50278.625
$ws.Range("N132").Value = -54229.178
$ws.Range("H134").Value = 1514.3889
$ws.Range("I134").Value = 1514.3889
$ws.Range("K134").Value = 4543.1667
$ws.Range("M134").Value = 526.8333000000002
$ws.Range("H136").Value = 3272
$ws.Range("I136").Value = 3272
$ws.Range("K136").Value = 9816
$ws.Range("M136").Value = -4716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26320880
$ws.Range("I70").Value = 38465570
$ws.Range("J70").Value = 7390.5
$ws.Range("K70").Value = 38465570
$ws.Range("L70").Value = 7390.5
$ws.Range("M70").Value = -38465300
$ws.Range("N70").Value = -7930.5
$ws.Range("H73").Value = 26320880
$ws.Range("I73").Value = 38465570
$ws.Range("J73").Value = 7390.5
$ws.Range("K73").Value = 38465570
$ws.Range("L73").Value = 7390.5
$ws.Range("M73").Value = -38464634
$ws.Range("N73").Value = -9262.5
$ws.Range("H122").Value = 47623430
$ws.Range("I122").Value = 76925260
$ws.Range("J122").Value = 7954.625
$ws.Range("K122").Value = 230775780
$ws.Range("L122").Value = 23863.875
$ws.Range("M122").Value = -230773330
$ws.Range("N122").Value = -28763.875
$ws.Range("H126").Value = 5308.8184
$ws.Range("I126").Value = 4499.625
$ws.Range("K126").Value = 13498.875
$ws.Range("M126").Value = -11028.875
$ws.Range("H132").Value = 2660.5454
$ws.Range("I132").Value = 2694.2273
$ws.Range("J132").Value = 2593.182
$ws.Range("K132").Value = 8082.6819
$ws.Range("L132").Value = 7779.545999999999
$ws.Range("M132").Value = -5552.6819
$ws.Range("N132").Value = -12839.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2499.5
$ws.Range("I7").Value = 2499.5
$ws.Range("K7").Value = 2499.5
$ws.Range("M7").Value = -2387.5
$ws.Range("H16").Value = 1717.75
$ws.Range("I16").Value = 1717.75
$ws.Range("K16").Value = 1717.75
$ws.Range("M16").Value = -1547.75
$ws.Range("H40").Value = 61441.906
$ws.Range("I40").Value = 83385.39999999999
$ws.Range("K40").Value = 83385.39999999999
$ws.Range("M40").Value = -83249.39999999999
$ws.Range("H126").Value = 2499.5
$ws.Range("I126").Value = 2499.5
$ws.Range("K126").Value = 7498.5
$ws.Range("M126").Value = -5028.5
$ws.Range("H128").Value = 45000
$ws.Range("I128").Value = 45000
$ws.Range("K128").Value = 45000
$ws.Range("M128").Value = -40020
$ws.Range("H132").Value = 7726.357
$ws.Range("I132").Value = 9040
$ws.Range("J132").Value = 6412.7144
$ws.Range("K132").Value = 27120
$ws.Range("L132").Value = 19238.1432
$ws.Range("M132").Value = -24590
$ws.Range("N132").Value = -24298.1432
$ws.Range("H136").Value = 6799.5
$ws.Range("I136").Value = 4700.6665
$ws.Range("K136").Value = 14101.9995
$ws.Range("M136").Value = -11551.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 43479620
$ws.Range("I100").Value = 1105.2
$ws.Range("K100").Value = 2210.4
$ws.Range("M100").Value = -1669.4
$ws.Range("H126").Value = 42459.168
$ws.Range("I126").Value = 61187.5
$ws.Range("J126").Value = 5002.5
$ws.Range("K126").Value = 183562.5
$ws.Range("L126").Value = 15007.5
$ws.Range("M126").Value = -181092.5
$ws.Range("N126").Value = -19947.5
$ws.Range("H132").Value = 5371.533
$ws.Range("I132").Value = 5213.3076
$ws.Range("K132").Value = 15639.9228
$ws.Range("M132").Value = -13109.9228
$ws.Range("H136").Value = 16927.8
$ws.Range("I136").Value = 18976.428
$ws.Range("K136").Value = 56929.284
$ws.Range("M136").Value = -54379.284
